# Apply numeric updates to the F column (attendee/view counts) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 54
$ws.Range("F4").Value = 7926
$ws.Range("F5").Value = 7926
$ws.Range("F6").Value = 103
$ws.Range("F8").Value = 2147
$ws.Range("F9").Value = 8672
$ws.Range("F10").Value = 13
$ws.Range("F13").Value = 5801
$ws.Range("F14").Value = 66
$ws.Range("F15").Value = 2801
$ws.Range("F16").Value = 1214
$ws.Range("F19").Value = 47
$ws.Range("F20").Value = 634
$ws.Range("F21").Value = 111
$ws.Range("F22").Value = 3989
$ws.Range("F23").Value = 84
$ws.Range("F25").Value = 69
$ws.Range("F27").Value = 182
$ws.Range("F28").Value = 34
$ws.Range("F29").Value = 5739
$ws.Range("F31").Value = 77
$ws.Range("F32").Value = 402
$ws.Range("F33").Value = 167
$ws.Range("F34").Value = 416
$ws.Range("F35").Value = 3196
$ws.Range("F36").Value = 1550
$ws.Range("F39").Value = 5699
$ws.Range("F42").Value = 55
$ws.Range("F43").Value = 3686
$ws.Range("F44").Value = 33
$ws.Range("F49").Value = 741
$ws.Range("F50").Value = 32

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 163
$ws.Range("F5").Value = 78
$ws.Range("F6").Value = 22
$ws.Range("F11").Value = 9

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1374

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1374
$ws.Range("F4").Value = 7926
$ws.Range("F5").Value = 7926
$ws.Range("F6").Value = 103
$ws.Range("F8").Value = 2147
$ws.Range("F9").Value = 8672
$ws.Range("F10").Value = 13
$ws.Range("F13").Value = 5801
$ws.Range("F14").Value = 66
$ws.Range("F15").Value = 2801
$ws.Range("F16").Value = 1214
$ws.Range("F19").Value = 47
$ws.Range("F20").Value = 163
$ws.Range("F21").Value = 634
$ws.Range("F22").Value = 111
$ws.Range("F23").Value = 3989
$ws.Range("F24").Value = 84
$ws.Range("F27").Value = 182
$ws.Range("F28").Value = 34
$ws.Range("F29").Value = 5739
$ws.Range("F30").Value = 77
$ws.Range("F31").Value = 402
$ws.Range("F32").Value = 167
$ws.Range("F33").Value = 416
$ws.Range("F34").Value = 78
$ws.Range("F35").Value = 3196
$ws.Range("F36").Value = 1550
$ws.Range("F37").Value = 22
$ws.Range("F41").Value = 5699
$ws.Range("F44").Value = 3686
$ws.Range("F49").Value = 745
$ws.Range("F50").Value = 9
